$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) updates: force text format to avoid numeric auto-conversion
# e.g. "620.59" must remain the literal text "620.59", not the number 620.59
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.952.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.656.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "620.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.655.32"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.540"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.69"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.504"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.268.06"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.677.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.965.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "520.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.745"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.63"
$ws.Range("D32").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "45.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.120.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "425.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "28.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.63"
$ws.Range("D49").Style = "Normal"

# Other text columns (Coin name, Link, Volume(1h)) updates
$ws.Range("E2").Value = "  +6.09%  "
$ws.Range("E3").Value = "  +17.86%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  +7.65%  "
$ws.Range("E6").Value = "  +2.25%  "
$ws.Range("E7").Value = "  +17.85%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  +5.03%  "
$ws.Range("E10").Value = "  +8.28%  "
$ws.Range("E11").Value = "  +5.32%  "
$ws.Range("E12").Value = "  +7.68%  "
$ws.Range("E13").Value = "  +11.70%  "
$ws.Range("E14").Value = "  +5.96%  "
$ws.Range("E15").Value = "  +17.90%  "
$ws.Range("E16").Value = "  +18.52%  "
$ws.Range("E17").Value = "  +6.06%  "
$ws.Range("E18").Value = "  +1.92%  "
$ws.Range("E20").Value = "  +8.21%  "
$ws.Range("E21").Value = "  +1.35%  "
$ws.Range("E22").Value = "  +18.53%  "
$ws.Range("E23").Value = "  +7.92%  "
$ws.Range("E24").Value = "  +5.97%  "
$ws.Range("E25").Value = "  +11.61%  "
$ws.Range("E26").Value = "  +6.99%  "
$ws.Range("E27").Value = "  +9.72%  "
$ws.Range("E29").Value = "  +12.14%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("E30").Value = "  +3.10%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("E31").Value = "  +12.08%  "
$ws.Range("E32").Value = "  +13.04%  "
$ws.Range("E33").Value = "  +17.35%  "
$ws.Range("E34").Value = "  +3.75%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  +9.48%  "
$ws.Range("E37").Value = "  +8.70%  "
$ws.Range("E38").Value = "  +10.91%  "
$ws.Range("E39").Value = "  +10.28%  "
$ws.Range("E40").Value = "  +6.92%  "
$ws.Range("E41").Value = "  +5.38%  "
$ws.Range("E42").Value = "  -5.00%  "
$ws.Range("E43").Value = "  +6.01%  "
$ws.Range("E44").Value = "  +11.41%  "
$ws.Range("E45").Value = "  +13.68%  "
$ws.Range("E46").Value = "  +3.80%  "
$ws.Range("E47").Value = "  +8.26%  "
$ws.Range("E48").Value = "  +12.10%  "
$ws.Range("E49").Value = "  +4.00%  "
$ws.Range("E51").Value = "  +10.17%  "
